$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting of the existing header cell (H1) onto the new ones
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the two new columns, row by row
$values = @{
  2  = @(7, 8)
  3  = @(1, 2)
  4  = @(7, 7)
  5  = @(8, 9)
  6  = @(7, 8)
  7  = @(8, 8)
  8  = @(9, 9)
  9  = @(7, 8)
  10 = @(9, 9)
  11 = @(8, 9)
  12 = @(7, 8)
  13 = @(6, 7)
  14 = @(5, 5)
  15 = @(1, 3)
  16 = @(4, 5)
}

foreach ($row in $values.Keys) {
  $pair = $values[$row]
  $ws.Cells.Item($row, 9).Value = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
}
